$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new "context" (H) and "type" (I) columns ---
$ws.Range("H1").Value = "context"
$ws.Range("H1").NumberFormat = "0.0000"
$ws.Range("I1").Value = "type"
$ws.Range("I1").NumberFormat = "0.0000"

$ws.Range("H2").Value = "Hadjuabudllah l1"
$ws.Range("H2").Font.Color = 0
$ws.Range("I2").Value = "soil layer 1"

$ws.Range("H3").Value = "Hadjuabudllah l2"
$ws.Range("H3").Font.Color = 0
$ws.Range("I3").Value = "soil layer 2"
$ws.Range("I3").Font.Color = 0

$ws.Range("H4").Value = "Laona soil l1"
$ws.Range("H4").Font.Color = 0
$ws.Range("I4").Value = "soil layer 1"
$ws.Range("I4").Font.Color = 0

$ws.Range("H5").Value = "Laona soil l2"
$ws.Range("H5").Font.Color = 0
$ws.Range("I5").Value = "soil layer 2"
$ws.Range("I5").Font.Color = 0

$ws.Range("H6").Value = "Laona soil l3"
$ws.Range("H6").Font.Color = 0
$ws.Range("I6").Value = "soil layer 3"
$ws.Range("I6").Font.Color = 0

$ws.Range("H7").Value = "Laona soil l1"
$ws.Range("H7").Font.Color = 0
$ws.Range("I7").Value = "soil layer 1"
$ws.Range("I7").Font.Color = 0

$ws.Range("H8").Value = "Laona soil l2"
$ws.Range("H8").Font.Color = 0
$ws.Range("I8").Value = "soil layer 2"
$ws.Range("I8").Font.Color = 0

$ws.Range("H9").Value = "Laona soil l3"
$ws.Range("H9").Font.Color = 0
$ws.Range("I9").Value = "soil layer 3"
$ws.Range("I9").Font.Color = 0

$ws.Range("H10").Value = "LA54:4"
$ws.Range("H10").Font.Color = 0
$ws.Range("I10").Value = "mudbrick"
$ws.Range("I10").HorizontalAlignment = -4131

$ws.Range("H11").Value = "LA54:4"
$ws.Range("H11").Font.Color = 0
$ws.Range("I11").Value = "mudbrick"
$ws.Range("I11").HorizontalAlignment = -4131

$ws.Range("H12").Value = "NA"
$ws.Range("H12").NumberFormat = "0.0000"
$ws.Range("I12").Value = "NA"
$ws.Range("I12").NumberFormat = "0.0000"

$ws.Range("H13").Value = "LA54:4"
$ws.Range("H13").Font.Color = 0
$ws.Range("I13").Value = "mudbrick"
$ws.Range("I13").HorizontalAlignment = -4131

$ws.Range("H14").Value = "LA54:4"
$ws.Range("H14").Font.Color = 0
$ws.Range("I14").Value = "mudbrick"
$ws.Range("I14").HorizontalAlignment = -4131

$ws.Range("H15").Value = "LA54:4"
$ws.Range("H15").Font.Color = 0
$ws.Range("I15").Value = "mudbrick"
$ws.Range("I15").HorizontalAlignment = -4131

$ws.Range("H16").Value = "LA54:4"
$ws.Range("H16").Font.Color = 0
$ws.Range("I16").Value = "mudbrick"
$ws.Range("I16").HorizontalAlignment = -4131

$ws.Range("H17").Value = "LA59:2"
$ws.Range("H17").Font.Color = 0
$ws.Range("I17").Value = "mudbrick"
$ws.Range("I17").HorizontalAlignment = -4131

$ws.Range("H18").Value = "LA59:2"
$ws.Range("H18").Font.Color = 0
$ws.Range("I18").Value = "mudbrick"
$ws.Range("I18").HorizontalAlignment = -4131

$ws.Range("H19").Value = "LA59:2"
$ws.Range("H19").Font.Color = 0
$ws.Range("I19").Value = "mudbrick"
$ws.Range("I19").HorizontalAlignment = -4131

$ws.Range("H20").Value = "LA59:2"
$ws.Range("H20").Font.Color = 0
$ws.Range("I20").Value = "mudbrick"
$ws.Range("I20").HorizontalAlignment = -4131

$ws.Range("H21").Value = "LA59:2"
$ws.Range("H21").Font.Color = 0
$ws.Range("I21").Value = "mudbrick"
$ws.Range("I21").HorizontalAlignment = -4131

$ws.Range("H22").Value = "LA54:7"
$ws.Range("H22").Font.Color = 0
$ws.Range("I22").Value = "mudbrick"
$ws.Range("I22").HorizontalAlignment = -4131

$ws.Range("H23").Value = "Hadjuabudllah l1"
$ws.Range("H23").Font.Color = 0
$ws.Range("I23").Value = "soil layer 1"

$ws.Range("H24").Value = "Hadjuabudllah l2"
$ws.Range("H24").Font.Color = 0
$ws.Range("I24").Value = "soil layer 2"
$ws.Range("I24").Font.Color = 0

$ws.Range("H25").Value = "Laona soil l1"
$ws.Range("H25").Font.Color = 0
$ws.Range("I25").Value = "soil layer 1"
$ws.Range("I25").Font.Color = 0

$ws.Range("H26").Value = "Laona soil l2"
$ws.Range("H26").Font.Color = 0
$ws.Range("I26").Value = "soil layer 2"
$ws.Range("I26").Font.Color = 0

$ws.Range("H27").Value = "Laona soil l3"
$ws.Range("H27").Font.Color = 0
$ws.Range("I27").Value = "soil layer 3"
$ws.Range("I27").Font.Color = 0

$ws.Range("H28").Value = "Laona soil l1"
$ws.Range("H28").Font.Color = 0
$ws.Range("I28").Value = "soil layer 1"
$ws.Range("I28").Font.Color = 0

$ws.Range("H29").Value = "Laona soil l2"
$ws.Range("H29").Font.Color = 0
$ws.Range("I29").Value = "soil layer 2"
$ws.Range("I29").Font.Color = 0

$ws.Range("H30").Value = "Laona soil l3"
$ws.Range("H30").Font.Color = 0
$ws.Range("I30").Value = "soil layer 3"
$ws.Range("I30").Font.Color = 0

$ws.Range("H31").Value = "NA"
$ws.Range("H31").NumberFormat = "0.0000"
$ws.Range("I31").Value = "NA"
$ws.Range("I31").NumberFormat = "0.0000"

# --- Update selection to match the saved view ---
$ws.Range("M11").Select()

Write-Output "done"
